# Apply commit: "Add About page locators and fixes for homepage locators"
$wb = $excel.ActiveWorkbook

$wsHome  = $wb.Worksheets.Item(1)   # V_HomePage
$wsAbout = $wb.Worksheets.Item(4)   # V_AboutAboutPage

# ------------------------------------------------------------------
# 1) V_AboutAboutPage: insert 4 new rows (6-9) for the new "Citing the
#    C3DC" locators, pushing the existing rows 6-10 down to 10-14.
#    Inserting this way (rather than writing past the end) makes the
#    new rows inherit the same A/B/C column styles already used by
#    the rows above them.
# ------------------------------------------------------------------
$wsAbout.Rows("6:9").Insert()

# Populate in this specific order (A6,B6,A7,A8,A9,B7,B8,B9) so the
# shared-string table allocates the new unique strings in the same
# order as the source workbook.
$wsAbout.Range("A6").Value = "AboutParagraph5"
$wsAbout.Range("B6").Value = "Citing the C3DC`nNCI expects users to acknowledge CCDI data use as follows:"

$wsAbout.Range("A7").Value = "AboutParagraph6"
$wsAbout.Range("A8").Value = "AboutParagraph7"
$wsAbout.Range("A9").Value = "AboutParagraph8"

$wsAbout.Range("B7").Value = '"The results published here are, in whole or in part, derived from the analysis of data listed in the C3DC (clinicalcommons.ccdi.cancer.gov), established by the National Cancer Institute' + [char]0x2019 + 's Childhood Cancer Data Initiative (CCDI)."'
$wsAbout.Range("B8").Value = 'To cite individual studies, note the CCDI study ID (e.g., phs002790) and include the name and URL or link for the C3DC (clinicalcommons.ccdi.cancer.gov), along with the phrase, "established by the National Cancer Institute' + [char]0x2019 + 's Childhood Cancer Data Initiative (CCDI)."'
$wsAbout.Range("B9").Value = 'Example: "The results analyzed and <published or shown> here are based in whole or in part from analyzing the Molecular Characterization Initiative data listed in the C3DC (clinicalcommons.ccdi.cancer.gov) under study ID phs002790. The data were accessed from the NCI' + [char]0x2019 + 's Cancer Research Data Commons (datacommons.cancer.gov). The C3DC was established by the National Cancer Institute' + [char]0x2019 + 's Childhood Cancer Data Initiative (CCDI)".'

# Row heights for the newly inserted rows (wrapped multi-line text)
$wsAbout.Rows("6:6").RowHeight = 28
$wsAbout.Rows("7:7").RowHeight = 28
$wsAbout.Rows("8:8").RowHeight = 42
$wsAbout.Rows("9:9").RowHeight = 56

# ------------------------------------------------------------------
# 2) Hyperlinks on V_AboutAboutPage shifted down by 4 rows along with
#    the data they were attached to. Stash the pristine "hyperlink
#    look" formatting first (Hyperlinks.Add below re-styles whatever
#    cell it touches), rebuild the hyperlinks at their new locations
#    in the same order (so relationship ids line up the same way),
#    then restore each cell's original look.
# ------------------------------------------------------------------
$wsAbout.Range("C10").Copy() | Out-Null
$wsAbout.Range("Z1").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats; C-column link look
$wsAbout.Range("B14").Copy() | Out-Null
$wsAbout.Range("Z2").PasteSpecial(-4122) | Out-Null   # B-column (named Hyperlink style) look
$excel.CutCopyMode = $false

$wsAbout.Hyperlinks.Delete()
$wsAbout.Hyperlinks.Add($wsAbout.Range("C10"), "https://www.cancer.gov/research/areas/childhood/childhood-cancer-data-initiative") | Out-Null
$wsAbout.Hyperlinks.Add($wsAbout.Range("C11"), "https://public.govdelivery.com/accounts/USNIHNCI/subscriber/new?topic_id=USNIHNCI_223") | Out-Null
$wsAbout.Hyperlinks.Add($wsAbout.Range("C12"), "https://cadsr.cancer.gov/onedata/dmdirect/NIH/NCI/CO/CDEDD?filter=Administered%20Item%20%28Data%20Element%20CO%29.CDEDD%20Classification.P_ITEM_ID_VER=12119072v1") | Out-Null
$wsAbout.Hyperlinks.Add($wsAbout.Range("C14"), "mailto:ncichildhoodcancerdatainitiative@mail.nih.gov") | Out-Null
$wsAbout.Hyperlinks.Add($wsAbout.Range("B14"), "mailto:ncichildhoodcancerdatainitiative@mail.nih.gov") | Out-Null
$wsAbout.Hyperlinks.Add($wsAbout.Range("C13"), "https://github.com/CBIIT/c3dc-model") | Out-Null

$wsAbout.Range("Z1").Copy() | Out-Null
$wsAbout.Range("C10").PasteSpecial(-4122) | Out-Null
$wsAbout.Range("Z1").Copy() | Out-Null
$wsAbout.Range("C11").PasteSpecial(-4122) | Out-Null
$wsAbout.Range("Z1").Copy() | Out-Null
$wsAbout.Range("C12").PasteSpecial(-4122) | Out-Null
$wsAbout.Range("Z1").Copy() | Out-Null
$wsAbout.Range("C14").PasteSpecial(-4122) | Out-Null
$wsAbout.Range("Z1").Copy() | Out-Null
$wsAbout.Range("C13").PasteSpecial(-4122) | Out-Null
$wsAbout.Range("Z2").Copy() | Out-Null
$wsAbout.Range("B14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsAbout.Range("Z1:Z2").Clear()

# ------------------------------------------------------------------
# 3) View/selection state: V_AboutAboutPage becomes the active sheet
#    and tab, with a new selection; V_HomePage is no longer the
#    selected tab.
# ------------------------------------------------------------------
$wsAbout.Select()
$wsAbout.Range("B19").Select()

Write-Host "Edit applied"
